$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.848.11"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "1.813.40"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.20"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4636"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3696"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07348"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8686"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.41"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "1.876.09"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.340"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07080"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.504"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.37"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008716"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.67"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("D21").Value = "26.906.76"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.337"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.55"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "2.118.03"
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.902"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.92"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.35"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.124"
$ws.Range("E28").Value = "  -5.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.299"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.34"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08896"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7566"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.149"
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.930"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.453"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.093"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01953"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05253"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.382"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5329"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.218"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1660"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.432"
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4934"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.29"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.670"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.01"
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06272"
$ws.Range("E51").Value = "  -0.81%  "
